$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 27.03890566666666
$ws.Range("H2").Value = 81.116717
$ws.Range("I2").Value = 0.07096188219033728
$ws.Range("J2").Value = 0.07096188219033729
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.261437
$ws.Range("N2").Value = 63.784311
$ws.Range("O2").Value = 0.8363123822876132
$ws.Range("P2").Value = 0.8363123822876132
$ws.Range("Q2").Value = 574.8859893807763
$ws.Range("R2").Value = 5173.973904426987
$ws.Range("S2").Value = 0.05934630074621392
$ws.Range("T2").Value = 0.05934630074621393

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 27.03890566666666
$ws.Range("H3").Value = 81.116717
$ws.Range("I3").Value = 0.07096188219033728
$ws.Range("J3").Value = 0.07096188219033729
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.321607
$ws.Range("N3").Value = 9.964821
$ws.Range("O3").Value = 0.1306544361603222
$ws.Range("P3").Value = 0.1306544361603222
$ws.Range("Q3").Value = 89.81261833473967
$ws.Range("R3").Value = 808.313565012657
$ws.Range("S3").Value = 0.009271484706453725
$ws.Range("T3").Value = 0.009271484706453727

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 27.03890566666666
$ws.Range("H4").Value = 81.116717
$ws.Range("I4").Value = 0.07096188219033728
$ws.Range("J4").Value = 0.07096188219033729
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.8397973333333333
$ws.Range("N4").Value = 2.519392
$ws.Range("O4").Value = 0.03303318155206465
$ws.Range("P4").Value = 0.03303318155206465
$ws.Range("Q4").Value = 22.70720087511822
$ws.Range("R4").Value = 204.364807876064
$ws.Range("S4").Value = 0.002344096737669634
$ws.Range("T4").Value = 0.002344096737669635

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 345.566579
$ws.Range("H5").Value = 1036.699737
$ws.Range("I5").Value = 0.9069174311350353
$ws.Range("J5").Value = 0.9069174311350354
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.261437
$ws.Range("N5").Value = 63.784311
$ws.Range("O5").Value = 0.8363123822876132
$ws.Range("P5").Value = 0.8363123822876132
$ws.Range("Q5").Value = 7347.242048714023
$ws.Range("R5").Value = 66125.1784384262
$ws.Range("S5").Value = 0.7584662773707037
$ws.Range("T5").Value = 0.7584662773707038

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 345.566579
$ws.Range("H6").Value = 1036.699737
$ws.Range("I6").Value = 0.9069174311350353
$ws.Range("J6").Value = 0.9069174311350354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.321607
$ws.Range("N6").Value = 9.964821
$ws.Range("O6").Value = 0.1306544361603222
$ws.Range("P6").Value = 0.1306544361603222
$ws.Range("Q6").Value = 1147.836367772453
$ws.Range("R6").Value = 10330.52730995208
$ws.Range("S6").Value = 0.1184927856089159
$ws.Range("T6").Value = 0.1184927856089159

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 345.566579
$ws.Range("H7").Value = 1036.699737
$ws.Range("I7").Value = 0.9069174311350353
$ws.Range("J7").Value = 0.9069174311350354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.8397973333333333
$ws.Range("N7").Value = 2.519392
$ws.Range("O7").Value = 0.03303318155206465
$ws.Range("P7").Value = 0.03303318155206465
$ws.Range("Q7").Value = 290.2058915333226
$ws.Range("R7").Value = 2611.853023799903
$ws.Range("S7").Value = 0.02995836815541571
$ws.Range("T7").Value = 0.02995836815541571

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.428738666666668
$ws.Range("H8").Value = 25.286216
$ws.Range("I8").Value = 0.0221206866746274
$ws.Range("J8").Value = 0.02212068667462741
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.261437
$ws.Range("N8").Value = 63.784311
$ws.Range("O8").Value = 0.8363123822876132
$ws.Range("P8").Value = 0.8363123822876132
$ws.Range("Q8").Value = 179.2070961507974
$ws.Range("R8").Value = 1612.863865357176
$ws.Range("S8").Value = 0.0184998041706955
$ws.Range("T8").Value = 0.01849980417069551

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.428738666666668
$ws.Range("H9").Value = 25.286216
$ws.Range("I9").Value = 0.0221206866746274
$ws.Range("J9").Value = 0.02212068667462741
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.321607
$ws.Range("N9").Value = 9.964821
$ws.Range("O9").Value = 0.1306544361603222
$ws.Range("P9").Value = 0.1306544361603222
$ws.Range("Q9").Value = 27.99695735637067
$ws.Range("R9").Value = 251.972616207336
$ws.Range("S9").Value = 0.002890165844952595
$ws.Range("T9").Value = 0.002890165844952596

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.428738666666668
$ws.Range("H10").Value = 25.286216
$ws.Range("I10").Value = 0.0221206866746274
$ws.Range("J10").Value = 0.02212068667462741
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.8397973333333333
$ws.Range("N10").Value = 2.519392
$ws.Range("O10").Value = 0.03303318155206465
$ws.Range("P10").Value = 0.03303318155206465
$ws.Range("Q10").Value = 7.078432255630223
$ws.Range("R10").Value = 63.705890300672
$ws.Range("S10").Value = 0.000730716658979304
$ws.Range("T10").Value = 0.0007307166589793042
